$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 78: 45205 -> 45206
$ws.Range("C2:C78").Value = 45206
